$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string need a text-format
# guard so Excel COM does not auto-convert them to Number (the source data
# keeps these columns as text/inlineStr even when the text looks numeric).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '63.452.32'
$ws.Range("E2").Value = '  +2.40%  '

$ws.Range("D3").Value = '3.469.84'
$ws.Range("E3").Value = '  +1.14%  '

$ws.Range("E4").Value = '  -0.07%  '

Set-TextValue $ws.Range("D5") '414.49'
$ws.Range("E5").Value = '  +1.00%  '

Set-TextValue $ws.Range("D6") '128.84'
$ws.Range("E6").Value = '  -1.01%  '

Set-TextValue $ws.Range("D7") '0.633'
$ws.Range("E7").Value = '  -0.26%  '

Set-TextValue $ws.Range("D8") '1.00'
$ws.Range("E8").Value = '  -0.01%  '

Set-TextValue $ws.Range("D9") '0.755'
$ws.Range("E9").Value = '  +2.34%  '

Set-TextValue $ws.Range("D10") '0.155'
$ws.Range("E10").Value = '  +10.79%  '

Set-TextValue $ws.Range("D11") '42.29'
$ws.Range("E11").Value = '  -3.22%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D12") '9.69'
$ws.Range("E12").Value = '  +3.51%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D13") '0.0000227'
$ws.Range("E13").Value = '  +1.33%  '

$ws.Range("D14").Value = '4.022.35'
$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("E15").Value = '  -1.04%  '

Set-TextValue $ws.Range("D16") '20.22'
$ws.Range("E16").Value = '  -4.76%  '

$ws.Range("D17").Value = '3.455.50'
$ws.Range("E17").Value = '  +1.21%  '

$ws.Range("E18").Value = '  +0.53%  '

Set-TextValue $ws.Range("D19") '12.32'
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("D20").Value = '63.406.33'
$ws.Range("E20").Value = '  +2.35%  '

Set-TextValue $ws.Range("D21") '454.74'
$ws.Range("E21").Value = '  -10.83%  '

Set-TextValue $ws.Range("D22") '89.67'
$ws.Range("E22").Value = '  -3.50%  '

Set-TextValue $ws.Range("D23") '3.27'
$ws.Range("E23").Value = '  -1.43%  '

Set-TextValue $ws.Range("D24") '13.11'
$ws.Range("E24").Value = '  -2.40%  '

Set-TextValue $ws.Range("D25") '10.09'
$ws.Range("E25").Value = '  +8.81%  '

Set-TextValue $ws.Range("D26") '3.29'
$ws.Range("E26").Value = '  -1.76%  '

Set-TextValue $ws.Range("D27") '33.38'
$ws.Range("E27").Value = '  -5.01%  '

Set-TextValue $ws.Range("D28") '4.77'
$ws.Range("E28").Value = '  -0.61%  '

Set-TextValue $ws.Range("D29") '12.44'
$ws.Range("E29").Value = '  +2.30%  '

$ws.Range("E30").Value = '  -2.10%  '

$ws.Range("E31").Value = '  -0.69%  '

Set-TextValue $ws.Range("D32") '0.168'
$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("E33").Value = '  -2.52%  '

Set-TextValue $ws.Range("D34") '39.89'
$ws.Range("E34").Value = '  -5.00%  '

$ws.Range("E35").Value = '  +0.00%  '

Set-TextValue $ws.Range("D36") '57.57'
$ws.Range("E36").Value = '  -2.77%  '

$ws.Range("E37").Value = '  -3.02%  '

Set-TextValue $ws.Range("D38") '3.08'
$ws.Range("E38").Value = '  +4.41%  '

Set-TextValue $ws.Range("D39") '1.00'
$ws.Range("E39").Value = '  +0.08%  '

Set-TextValue $ws.Range("D40") '2.80'
$ws.Range("E40").Value = '  +1.46%  '

Set-TextValue $ws.Range("D41") '0.136'
$ws.Range("E41").Value = '  -1.89%  '

Set-TextValue $ws.Range("D42") '4.52'
$ws.Range("E42").Value = '  +4.29%  '

$ws.Range("D43").Value = '0.0₃0646'
$ws.Range("E43").Value = '  +54.95%  '

Set-TextValue $ws.Range("D44") '146.20'
$ws.Range("E44").Value = '  -1.21%  '

Set-TextValue $ws.Range("D45") '3.31'
$ws.Range("E45").Value = '  -4.39%  '

Set-TextValue $ws.Range("D46") '0.314'
$ws.Range("E46").Value = '  -1.58%  '

Set-TextValue $ws.Range("D47") '1.99'
$ws.Range("E47").Value = '  -6.25%  '

Set-TextValue $ws.Range("D48") '2.33'
$ws.Range("E48").Value = '  -2.40%  '

Set-TextValue $ws.Range("D49") '15.98'
$ws.Range("E49").Value = '  -3.98%  '

Set-TextValue $ws.Range("D50") '21.48'
$ws.Range("E50").Value = '  -6.43%  '

Set-TextValue $ws.Range("D51") '0.139'
$ws.Range("E51").Value = '  -4.98%  '
